$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(45534.5, 0.3513, 0.3546, 0.3383, 0.34, 47841509.2),
  @(45534.66666666666, 0.34, 0.3493, 0.3373, 0.3455, 16638411.6),
  @(45534.83333333334, 0.3455, 0.3482, 0.3455, 0.3469, 4426300.1),
  @(45535, 0.3469, 0.3485, 0.3458, 0.3469, 1831291.6),
  @(45535.16666666666, 0.3474, 0.3504, 0.3474, 0.3478, 5161715.8),
  @(45535.33333333334, 0.3478, 0.3495, 0.3453, 0.3481, 4809648.4),
  @(45535.5, 0.3481, 0.3496, 0.3471, 0.3492, 4712130.7),
  @(45535.66666666666, 0.3492, 0.3493, 0.345, 0.3468, 4751588.4),
  @(45535.83333333334, 0.3467, 0.347, 0.3447, 0.345, 1160830.8),
  @(45536, 0.3451, 0.3458, 0.3394, 0.3448, 6403007.4),
  @(45536.16666666666, 0.3449, 0.3449, 0.3393, 0.3446, 5978305.4),
  @(45536.33333333334, 0.3446, 0.345, 0.3428, 0.3428, 1772825.1),
  @(45536.5, 0.3431, 0.3436, 0.3351, 0.3403, 15760883.4),
  @(45536.66666666666, 0.3403, 0.3422, 0.3341, 0.3397, 13694886.7),
  @(45536.83333333334, 0.3397, 0.3421, 0.3267, 0.3311, 17389891),
  @(45537, 0.331, 0.3345, 0.3308, 0.3316, 6640945.6),
  @(45537.16666666666, 0.3315, 0.3316, 0.325, 0.3259, 11315069.6),
  @(45537.33333333334, 0.3259, 0.3364, 0.3253, 0.3346, 19529681.1),
  @(45537.5, 0.3338, 0.3349, 0.3269, 0.3325, 15387478.4),
  @(45537.66666666666, 0.3325, 0.3342, 0.3297, 0.3309, 9306857.800000001),
  @(45537.83333333334, 0.3309, 0.3374, 0.3303, 0.3362, 8010292.1),
  @(45538, 0.3362, 0.3398, 0.3342, 0.3351, 10001431.3),
  @(45538.16666666666, 0.3352, 0.3356, 0.3318, 0.3329, 6379042.4),
  @(45538.33333333334, 0.3335, 0.3336, 0.3281, 0.3302, 10360581.5),
  @(45538.5, 0.3302, 0.332, 0.32, 0.3204, 26948259.2),
  @(45538.66666666666, 0.3204, 0.3259, 0.3195, 0.3249, 9479136.9),
  @(45538.83333333334, 0.3249, 0.3262, 0.3172, 0.3185, 11673959.8),
  @(45539, 0.3184, 0.3202, 0.305, 0.3157, 44782073.1),
  @(45539.16666666666, 0.3157, 0.3222, 0.3124, 0.3212, 16564924.7),
  @(45539.33333333334, 0.3212, 0.325, 0.315, 0.3175, 16756179.9),
  @(45539.5, 0.3176, 0.3262, 0.3154, 0.326, 17387851.2),
  @(45539.66666666666, 0.3259, 0.3293, 0.321, 0.3223, 14170383.3),
  @(45539.83333333334, 0.3231, 0.3234, 0.3219, 0.3223, 940070.4)
)

$startRow = 1085
$formatSourceRow = 1084

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Range("A$formatSourceRow`:F$formatSourceRow").Copy()
    $ws.Range("A$r`:F$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
}

$excel.CutCopyMode = 0
